$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 734.96
$ws.Range("I107").Value = 859.45
$ws.Range("J107").Value = 237
$ws.Range("K107").Value = 859.45
$ws.Range("L107").Value = 237
$ws.Range("M107").Value = 1060.55
$ws.Range("N107").Value = -4077
$ws.Range("H111").Value = 7145.0625
$ws.Range("I111").Value = 8061.875
$ws.Range("J111").Value = 6228.25
$ws.Range("K111").Value = 24185.625
$ws.Range("L111").Value = 18684.75
$ws.Range("M111").Value = -21118.625
$ws.Range("N111").Value = -24818.75
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = ""
$ws.Range("N121").Value = ""
$ws.Range("H137").Value = 2676.9211
$ws.Range("I137").Value = 2248.9355
$ws.Range("J137").Value = 4572.2856
$ws.Range("K137").Value = 6746.806500000001
$ws.Range("L137").Value = 13716.8568
$ws.Range("M137").Value = -4196.806500000001
$ws.Range("N137").Value = -18816.8568
$ws.Range("H138").Value = 2138.803
$ws.Range("I138").Value = 1479.6774
$ws.Range("J138").Value = 2722.6
$ws.Range("K138").Value = 4439.0322
$ws.Range("L138").Value = 8167.799999999999
$ws.Range("M138").Value = 700.9678000000004
$ws.Range("N138").Value = -18447.8
$ws.Range("H141").Value = 2965.087
$ws.Range("I141").Value = 842.3946999999999
$ws.Range("J141").Value = 13047.875
$ws.Range("K141").Value = 2527.1841
$ws.Range("L141").Value = 39143.625
$ws.Range("M141").Value = 2652.8159
$ws.Range("N141").Value = -49503.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2109.14
$ws.Range("I32").Value = 1741.2211
$ws.Range("J32").Value = 9099.6
$ws.Range("K32").Value = 1741.2211
$ws.Range("L32").Value = 9099.6
$ws.Range("M32").Value = -1454.2211
$ws.Range("N32").Value = -9673.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 454.5
$ws.Range("I5").Value = 404
$ws.Range("K5").Value = 404
$ws.Range("M5").Value = -291
$ws.Range("H94").Value = 1316.6666
$ws.Range("I94").Value = 1300
$ws.Range("J94").Value = 1325
$ws.Range("K94").Value = 1300
$ws.Range("L94").Value = 1325
$ws.Range("M94").Value = -849
$ws.Range("N94").Value = -2227
$ws.Range("H107").Value = 126486.375
$ws.Range("I107").Value = 251097.75
$ws.Range("J107").Value = 1875
$ws.Range("K107").Value = 251097.75
$ws.Range("L107").Value = 1875
$ws.Range("M107").Value = -249177.75
$ws.Range("N107").Value = -5715
$ws.Range("H134").Value = 2305.5881
$ws.Range("I134").Value = 1973.6086
$ws.Range("J134").Value = 2999.7273
$ws.Range("K134").Value = 5920.825800000001
$ws.Range("L134").Value = 8999.1819
$ws.Range("M134").Value = -3385.825800000001
$ws.Range("N134").Value = -14069.1819
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4720.4087
$ws.Range("I31").Value = 1174.5428
$ws.Range("J31").Value = 8167.778
$ws.Range("K31").Value = 1174.5428
$ws.Range("L31").Value = 8167.778
$ws.Range("M31").Value = -879.5427999999999
$ws.Range("N31").Value = -8757.778
$ws.Range("H34").Value = 4720.4087
$ws.Range("I34").Value = 1174.5428
$ws.Range("J34").Value = 8167.778
$ws.Range("K34").Value = 1174.5428
$ws.Range("L34").Value = 8167.778
$ws.Range("M34").Value = -972.5427999999999
$ws.Range("N34").Value = -8571.778
$ws.Range("H43").Value = 194000
$ws.Range("J43").Value = 194000
$ws.Range("L43").Value = 194000
$ws.Range("N43").Value = -194368
$ws.Range("H58").Value = 1416.5
$ws.Range("I58").Value = 1230.4375
$ws.Range("J58").Value = 1714.2
$ws.Range("K58").Value = 1230.4375
$ws.Range("L58").Value = 1714.2
$ws.Range("M58").Value = -1027.4375
$ws.Range("N58").Value = -2120.2
$ws.Range("H62").Value = 3888.889
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 3750
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 3750
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -4998
$ws.Range("H65").Value = 3888.889
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 3750
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 18750
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -24990
$ws.Range("H92").Value = 58299.5
$ws.Range("J92").Value = 58299.5
$ws.Range("L92").Value = 58299.5
$ws.Range("N92").Value = -63291.5
$ws.Range("H93").Value = 15069
$ws.Range("I93").Value = 6203.5
$ws.Range("J93").Value = 32800
$ws.Range("K93").Value = 6203.5
$ws.Range("L93").Value = 32800
$ws.Range("M93").Value = -4331.5
$ws.Range("N93").Value = -36544
$ws.Range("H95").Value = 20311.5
$ws.Range("J95").Value = 20311.5
$ws.Range("L95").Value = 20311.5
$ws.Range("N95").Value = -25803.5
$ws.Range("H96").Value = 45000
$ws.Range("J96").Value = 45000
$ws.Range("L96").Value = 45000
$ws.Range("N96").Value = -50492
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = ""
$ws.Range("H101").Value = 194000
$ws.Range("J101").Value = 194000
$ws.Range("L101").Value = 194000
$ws.Range("N101").Value = -200490
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""
$ws.Range("H103").Value = 100000
$ws.Range("J103").Value = 100000
$ws.Range("L103").Value = 100000
$ws.Range("M103").Value = -102344
$ws.Range("H105").Value = 2666.3333
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""
$ws.Range("H108").Value = 33066.5
$ws.Range("J108").Value = 33066.5
$ws.Range("L108").Value = 33066.5
$ws.Range("N108").Value = -40746.5
$ws.Range("H122").Value = 1829.9546
$ws.Range("I122").Value = 1790.45
$ws.Range("J122").Value = 1862.875
$ws.Range("K122").Value = 5371.35
$ws.Range("L122").Value = 5588.625
$ws.Range("M122").Value = -2921.35
$ws.Range("N122").Value = -10488.625
$ws.Range("H132").Value = 1414.925
$ws.Range("I132").Value = 987.25806
$ws.Range("J132").Value = 2888
$ws.Range("K132").Value = 2961.77418
$ws.Range("L132").Value = 8664
$ws.Range("M132").Value = -431.7741799999999
$ws.Range("N132").Value = -13724
$ws.Range("H136").Value = 1416.5
$ws.Range("I136").Value = 1230.4375
$ws.Range("J136").Value = 1714.2
$ws.Range("K136").Value = 3691.3125
$ws.Range("L136").Value = 5142.6
$ws.Range("M136").Value = -1141.3125
$ws.Range("N136").Value = -10242.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20002124
$ws.Range("J4").Value = 25002250
$ws.Range("L4").Value = 75006750
$ws.Range("N4").Value = -75006974
$ws.Range("H5").Value = 1144.5385
$ws.Range("I5").Value = 577.5
$ws.Range("K5").Value = 1732.5
$ws.Range("M5").Value = -1620.5
$ws.Range("H131").Value = 1130.5
$ws.Range("J131").Value = 1185.125
$ws.Range("L131").Value = 3555.375
$ws.Range("N131").Value = -13635.375
$ws.Range("H134").Value = 7143.6895
$ws.Range("I134").Value = 4843.3335
$ws.Range("J134").Value = 7743.7827
$ws.Range("K134").Value = 14530.0005
$ws.Range("L134").Value = 23231.3481
$ws.Range("M134").Value = -9460.000499999998
$ws.Range("N134").Value = -33371.3481
$ws.Range("H135").Value = 1144.5385
$ws.Range("I135").Value = 577.5
$ws.Range("K135").Value = 5197.5
$ws.Range("M135").Value = -2662.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2190.9092
$ws.Range("H132").Value = 2194.111
$ws.Range("I132").Value = 2164.476
$ws.Range("J132").Value = 2235.6
$ws.Range("K132").Value = 6493.428
$ws.Range("L132").Value = 6706.799999999999
$ws.Range("M132").Value = -3963.428
$ws.Range("N132").Value = -11766.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 8562
$ws.Range("I93").Value = 9713.666999999999
$ws.Range("J93").Value = 1652
$ws.Range("K93").Value = 9713.666999999999
$ws.Range("L93").Value = 1652
$ws.Range("M93").Value = -8465.666999999999
$ws.Range("N93").Value = -4148
$ws.Range("H101").Value = 30362
$ws.Range("J101").Value = 30362
$ws.Range("L101").Value = 30362
$ws.Range("N101").Value = -36852
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5026.2
$ws.Range("I96").Value = 3940
$ws.Range("J96").Value = 5226.2896
$ws.Range("K96").Value = 3940
$ws.Range("L96").Value = 5226.2896
$ws.Range("M96").Value = -2567
$ws.Range("N96").Value = -7972.2896
$ws.Range("H98").Value = 99750
$ws.Range("J98").Value = 99750
$ws.Range("L98").Value = 99750
$ws.Range("N98").Value = -105740
$ws.Range("H101").Value = 20050.75
$ws.Range("J101").Value = 20050.75
$ws.Range("L101").Value = 20050.75
$ws.Range("N101").Value = -26540.75
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").Value = ""
